$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.499.03'
$ws.Range('E2').Value = '  -2.42%  '
$ws.Range('D3').Value = '2.481.03'
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('E4').Value = '  +0.53%  '
$ws.Range('D5').Value = '''313.96'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').Value = '''93.23'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -5.71%  '
$ws.Range('D7').Value = '''0.547'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -2.73%  '
$ws.Range('D8').Value = '''1.01'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +0.46%  '
$ws.Range('D9').Value = '''0.495'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  -4.07%  '
$ws.Range('D10').Value = '''33.29'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -5.28%  '
$ws.Range('D11').Value = '''0.0781'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -2.29%  '
$ws.Range('D12').Value = '''0.109'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +0.12%  '
$ws.Range('D13').Value = '2.869.01'
$ws.Range('E13').Value = '  -1.40%  '
$ws.Range('D14').Value = '''6.89'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  -3.97%  '
$ws.Range('D15').Value = '''15.37'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').Value = '2.483.19'
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('D17').Value = '''0.787'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -2.47%  '
$ws.Range('D18').Value = '41.374.94'
$ws.Range('E18').Value = '  -2.72%  '
$ws.Range('D19').Value = '''6.30'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -4.47%  '
$ws.Range('D20').Value = '0.0₃0926'
$ws.Range('E20').Value = '  -1.28%  '
$ws.Range('D21').Value = '''70.09'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +1.63%  '
$ws.Range('D22').Value = '''11.14'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -8.39%  '
$ws.Range('D23').Value = '''234.90'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -2.65%  '
$ws.Range('D24').Value = '''2.75'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -3.51%  '
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').Value = '''1.89'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -5.12%  '
$ws.Range('D27').Value = '''24.12'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -5.21%  '
$ws.Range('D28').Value = '''2.25'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = '''9.79'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('D30').Value = '''36.56'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -3.13%  '
$ws.Range('D31').Value = '''153.40'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -2.17%  '
$ws.Range('D32').Value = '''5.47'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -6.66%  '
$ws.Range('D33').Value = '''2.56'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -2.82%  '
$ws.Range('D34').Value = '''0.0751'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -3.95%  '
$ws.Range('B35').Value = 'ApeXProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D35').Value = '''2.51'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -6.52%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').Value = '''17.78'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +1.01%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '''3.03'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -3.50%  '
$ws.Range('D38').Value = '''1.86'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -5.70%  '
$ws.Range('E39').Value = '  -3.28%  '
$ws.Range('D40').Value = '''0.100'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -7.00%  '
$ws.Range('D41').Value = '''4.07'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -2.49%  '
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('D43').Value = '''19.64'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -11.33%  '
$ws.Range('D44').Value = '1.976.39'
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('D45').Value = '''0.0283'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -4.15%  '
$ws.Range('D46').Value = '''2.97'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -7.56%  '
$ws.Range('D47').Value = '''8.79'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -1.80%  '
$ws.Range('D48').Value = '2.734.18'
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').Value = '''68.64'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -3.92%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '''96.24'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -3.58%  '
$ws.Range('D51').Value = '''0.177'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -5.92%  '
